# ============================================================================
# edit.ps1
#
# Commit: "feat: add 2022-Q1 data"
#
# The workbook's last sheet (6th, named "总计" / totals) is repurposed to
# hold the per-fund 2022-Q1 holdings detail (same shape as the 2021-Qx detail
# sheets) and gets renamed to "2022-Q1". A brand-new "总计" sheet is appended
# after it, containing the original totals table plus one new leading row
# summarizing the 2022-Q1 quarter.
# ============================================================================

function Set-TextCell($ws, $addr, $val) {
    # Force the value to be stored as TEXT. Left as a normal .Value
    # assignment, Excel auto-coerces numeric-looking strings (e.g. "25.17",
    # "050001") into numbers/doubles, which loses leading zeros and exact
    # decimal text. Temporarily marking the cell as Text ("@") keeps the
    # literal string, then ClearFormats() drops the temporary number format
    # again so the cell ends up with the workbook's default (unstyled) look.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

function Set-NumCell($ws, $addr, $val) {
    # Plain numeric cell, default style.
    $ws.Range($addr).Value = $val
}

function Set-StyledCell($ws, $addr, $val, $styleSrc) {
    # Cell (numeric index column, or header text) that must carry the same
    # bold/centered/bordered style used by the sheet's header row / index
    # column (cellXf "2" in the original workbook).
    $styleSrc.Copy()
    $dst = $ws.Range($addr)
    $dst.PasteSpecial(-4122)
    $dst.Value = $val
}

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Step 1 - the sheet currently named "总计" (position 6) holds the OLD
# totals table. Grab a style-source cell from it (A2: bold/centered/
# bordered header-like style) before overwriting its contents, then rename
# it to "2022-Q1" and replace its data with the fund-holdings detail table.
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
$styleSrc = $ws.Range("A2")
$ws.Name = "2022-Q1"

# -- header row --
$addr = "B1"; $val = "基金代码"
Set-StyledCell $ws $addr $val $styleSrc
$addr = "C1"; $val = "基金名称"
Set-StyledCell $ws $addr $val $styleSrc
$addr = "D1"; $val = "基金规模"
Set-StyledCell $ws $addr $val $styleSrc
$addr = "E1"; $val = "股票总仓位"
Set-StyledCell $ws $addr $val $styleSrc
$addr = "F1"; $val = "仓位占比"
Set-StyledCell $ws $addr $val $styleSrc
$addr = "G1"; $val = "持有市值(亿元)"
Set-StyledCell $ws $addr $val $styleSrc
$addr = "H1"; $val = "仓位排名"
Set-StyledCell $ws $addr $val $styleSrc

# -- data rows (fund holdings, 2022-Q1) --
# row 2 / index 0
$addr = "A2"; $val = 0
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B2"; $val = "516150"
Set-TextCell $ws $addr $val
$addr = "C2"; $val = "嘉实中证稀土产业ETF"
Set-TextCell $ws $addr $val
$addr = "D2"; $val = "25.17"
Set-TextCell $ws $addr $val
$addr = "E2"; $val = "99.75"
Set-TextCell $ws $addr $val
$addr = "F2"; $val = "6.47"
Set-TextCell $ws $addr $val
$addr = "G2"; $val = "1.6285"
Set-TextCell $ws $addr $val
$addr = "H2"; $val = 2
Set-NumCell $ws $addr $val
# row 3 / index 1
$addr = "A3"; $val = 1
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B3"; $val = "512400"
Set-TextCell $ws $addr $val
$addr = "C3"; $val = "南方中证申万有色金属ETF"
Set-TextCell $ws $addr $val
$addr = "D3"; $val = "36.45"
Set-TextCell $ws $addr $val
$addr = "E3"; $val = "99.71"
Set-TextCell $ws $addr $val
$addr = "F3"; $val = "3.50"
Set-TextCell $ws $addr $val
$addr = "G3"; $val = "1.2758"
Set-TextCell $ws $addr $val
$addr = "H3"; $val = 6
Set-NumCell $ws $addr $val
# row 4 / index 2
$addr = "A4"; $val = 2
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B4"; $val = "050001"
Set-TextCell $ws $addr $val
$addr = "C4"; $val = "博时价值增长混合"
Set-TextCell $ws $addr $val
$addr = "D4"; $val = "22.97"
Set-TextCell $ws $addr $val
$addr = "E4"; $val = "69.99"
Set-TextCell $ws $addr $val
$addr = "F4"; $val = "3.97"
Set-TextCell $ws $addr $val
$addr = "G4"; $val = "0.9119"
Set-TextCell $ws $addr $val
$addr = "H4"; $val = 1
Set-NumCell $ws $addr $val
# row 5 / index 3
$addr = "A5"; $val = 3
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B5"; $val = "160221"
Set-TextCell $ws $addr $val
$addr = "C5"; $val = "国泰国证有色金属行业指数（LOF）A"
Set-TextCell $ws $addr $val
$addr = "D5"; $val = "25.54"
Set-TextCell $ws $addr $val
$addr = "E5"; $val = "94.76"
Set-TextCell $ws $addr $val
$addr = "F5"; $val = "3.51"
Set-TextCell $ws $addr $val
$addr = "G5"; $val = "0.8965"
Set-TextCell $ws $addr $val
$addr = "H5"; $val = 5
Set-NumCell $ws $addr $val
# row 6 / index 4
$addr = "A6"; $val = 4
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B6"; $val = "014224"
Set-TextCell $ws $addr $val
$addr = "C6"; $val = "大成聚优成长混合A"
Set-TextCell $ws $addr $val
$addr = "D6"; $val = "33.72"
Set-TextCell $ws $addr $val
$addr = "E6"; $val = "53.38"
Set-TextCell $ws $addr $val
$addr = "F6"; $val = "2.47"
Set-TextCell $ws $addr $val
$addr = "G6"; $val = "0.8329"
Set-TextCell $ws $addr $val
$addr = "H6"; $val = 8
Set-NumCell $ws $addr $val
# row 7 / index 5
$addr = "A7"; $val = 5
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B7"; $val = "165520"
Set-TextCell $ws $addr $val
$addr = "C7"; $val = "中信保诚中证800 有色指数（LOF）"
Set-TextCell $ws $addr $val
$addr = "D7"; $val = "19.65"
Set-TextCell $ws $addr $val
$addr = "E7"; $val = "94.35"
Set-TextCell $ws $addr $val
$addr = "F7"; $val = "4.19"
Set-TextCell $ws $addr $val
$addr = "G7"; $val = "0.8233"
Set-TextCell $ws $addr $val
$addr = "H7"; $val = 6
Set-NumCell $ws $addr $val
# row 8 / index 6
$addr = "A8"; $val = 6
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B8"; $val = "516780"
Set-TextCell $ws $addr $val
$addr = "C8"; $val = "华泰柏瑞中证稀土产业ETF"
Set-TextCell $ws $addr $val
$addr = "D8"; $val = "11.06"
Set-TextCell $ws $addr $val
$addr = "E8"; $val = "98.70"
Set-TextCell $ws $addr $val
$addr = "F8"; $val = "6.39"
Set-TextCell $ws $addr $val
$addr = "G8"; $val = "0.7067"
Set-TextCell $ws $addr $val
$addr = "H8"; $val = 2
Set-NumCell $ws $addr $val
# row 9 / index 7
$addr = "A9"; $val = 7
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B9"; $val = "160526"
Set-TextCell $ws $addr $val
$addr = "C9"; $val = "博时优势企业3年封闭运作灵活配置混合（LOF）A"
Set-TextCell $ws $addr $val
$addr = "D9"; $val = "14.81"
Set-TextCell $ws $addr $val
$addr = "E9"; $val = "82.86"
Set-TextCell $ws $addr $val
$addr = "F9"; $val = "3.96"
Set-TextCell $ws $addr $val
$addr = "G9"; $val = "0.5865"
Set-TextCell $ws $addr $val
$addr = "H9"; $val = 9
Set-NumCell $ws $addr $val
# row 10 / index 8
$addr = "A10"; $val = 8
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B10"; $val = "050201"
Set-TextCell $ws $addr $val
$addr = "C10"; $val = "博时价值增长贰号混合"
Set-TextCell $ws $addr $val
$addr = "D10"; $val = "9.29"
Set-TextCell $ws $addr $val
$addr = "E10"; $val = "62.31"
Set-TextCell $ws $addr $val
$addr = "F10"; $val = "3.47"
Set-TextCell $ws $addr $val
$addr = "G10"; $val = "0.3224"
Set-TextCell $ws $addr $val
$addr = "H10"; $val = 1
Set-NumCell $ws $addr $val
# row 11 / index 9
$addr = "A11"; $val = 9
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B11"; $val = "159715"
Set-TextCell $ws $addr $val
$addr = "C11"; $val = "易方达中证稀土产业ETF"
Set-TextCell $ws $addr $val
$addr = "D11"; $val = "3.42"
Set-TextCell $ws $addr $val
$addr = "E11"; $val = "99.06"
Set-TextCell $ws $addr $val
$addr = "F11"; $val = "6.41"
Set-TextCell $ws $addr $val
$addr = "G11"; $val = "0.2192"
Set-TextCell $ws $addr $val
$addr = "H11"; $val = 2
Set-NumCell $ws $addr $val
# row 12 / index 10
$addr = "A12"; $val = 10
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B12"; $val = "159713"
Set-TextCell $ws $addr $val
$addr = "C12"; $val = "富国中证稀土产业交易型开放式指数证券投资基金"
Set-TextCell $ws $addr $val
$addr = "D12"; $val = "3.26"
Set-TextCell $ws $addr $val
$addr = "E12"; $val = "99.26"
Set-TextCell $ws $addr $val
$addr = "F12"; $val = "6.45"
Set-TextCell $ws $addr $val
$addr = "G12"; $val = "0.2103"
Set-TextCell $ws $addr $val
$addr = "H12"; $val = 2
Set-NumCell $ws $addr $val
# row 13 / index 11
$addr = "A13"; $val = 11
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B13"; $val = "510410"
Set-TextCell $ws $addr $val
$addr = "C13"; $val = "博时上证自然资源ETF"
Set-TextCell $ws $addr $val
$addr = "D13"; $val = "4.67"
Set-TextCell $ws $addr $val
$addr = "E13"; $val = "98.61"
Set-TextCell $ws $addr $val
$addr = "F13"; $val = "3.76"
Set-TextCell $ws $addr $val
$addr = "G13"; $val = "0.1756"
Set-TextCell $ws $addr $val
$addr = "H13"; $val = 6
Set-NumCell $ws $addr $val
# row 14 / index 12
$addr = "A14"; $val = 12
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B14"; $val = "011630"
Set-TextCell $ws $addr $val
$addr = "C14"; $val = "西藏东财中证有色金属指数增强A"
Set-TextCell $ws $addr $val
$addr = "D14"; $val = "3.87"
Set-TextCell $ws $addr $val
$addr = "E14"; $val = "92.28"
Set-TextCell $ws $addr $val
$addr = "F14"; $val = "4.52"
Set-TextCell $ws $addr $val
$addr = "G14"; $val = "0.1749"
Set-TextCell $ws $addr $val
$addr = "H14"; $val = 6
Set-NumCell $ws $addr $val
# row 15 / index 13
$addr = "A15"; $val = 13
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B15"; $val = "217001"
Set-TextCell $ws $addr $val
$addr = "C15"; $val = "招商安泰混合"
Set-TextCell $ws $addr $val
$addr = "D15"; $val = "4.22"
Set-TextCell $ws $addr $val
$addr = "E15"; $val = "70.79"
Set-TextCell $ws $addr $val
$addr = "F15"; $val = "3.45"
Set-TextCell $ws $addr $val
$addr = "G15"; $val = "0.1456"
Set-TextCell $ws $addr $val
$addr = "H15"; $val = 3
Set-NumCell $ws $addr $val
# row 16 / index 14
$addr = "A16"; $val = 14
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B16"; $val = "217012"
Set-TextCell $ws $addr $val
$addr = "C16"; $val = "招商行业领先混合A"
Set-TextCell $ws $addr $val
$addr = "D16"; $val = "2.86"
Set-TextCell $ws $addr $val
$addr = "E16"; $val = "80.77"
Set-TextCell $ws $addr $val
$addr = "F16"; $val = "4.59"
Set-TextCell $ws $addr $val
$addr = "G16"; $val = "0.1313"
Set-TextCell $ws $addr $val
$addr = "H16"; $val = 10
Set-NumCell $ws $addr $val
# row 17 / index 15
$addr = "A17"; $val = 15
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B17"; $val = "960019"
Set-TextCell $ws $addr $val
$addr = "C17"; $val = "招商行业领先混合H"
Set-TextCell $ws $addr $val
$addr = "D17"; $val = "2.86"
Set-TextCell $ws $addr $val
$addr = "E17"; $val = "80.77"
Set-TextCell $ws $addr $val
$addr = "F17"; $val = "4.59"
Set-TextCell $ws $addr $val
$addr = "G17"; $val = "0.1313"
Set-TextCell $ws $addr $val
$addr = "H17"; $val = 10
Set-NumCell $ws $addr $val
# row 18 / index 16
$addr = "A18"; $val = 16
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B18"; $val = "014225"
Set-TextCell $ws $addr $val
$addr = "C18"; $val = "大成聚优成长混合C"
Set-TextCell $ws $addr $val
$addr = "D18"; $val = "4.85"
Set-TextCell $ws $addr $val
$addr = "E18"; $val = "53.38"
Set-TextCell $ws $addr $val
$addr = "F18"; $val = "2.47"
Set-TextCell $ws $addr $val
$addr = "G18"; $val = "0.1198"
Set-TextCell $ws $addr $val
$addr = "H18"; $val = 8
Set-NumCell $ws $addr $val
# row 19 / index 17
$addr = "A19"; $val = 17
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B19"; $val = "011631"
Set-TextCell $ws $addr $val
$addr = "C19"; $val = "西藏东财中证有色金属指数增强C"
Set-TextCell $ws $addr $val
$addr = "D19"; $val = "2.18"
Set-TextCell $ws $addr $val
$addr = "E19"; $val = "92.28"
Set-TextCell $ws $addr $val
$addr = "F19"; $val = "4.52"
Set-TextCell $ws $addr $val
$addr = "G19"; $val = "0.0985"
Set-TextCell $ws $addr $val
$addr = "H19"; $val = 6
Set-NumCell $ws $addr $val
# row 20 / index 18
$addr = "A20"; $val = 18
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B20"; $val = "161217"
Set-TextCell $ws $addr $val
$addr = "C20"; $val = "国投瑞银中证上游资源产业指数(LOF)"
Set-TextCell $ws $addr $val
$addr = "D20"; $val = "3.10"
Set-TextCell $ws $addr $val
$addr = "E20"; $val = "94.10"
Set-TextCell $ws $addr $val
$addr = "F20"; $val = "2.56"
Set-TextCell $ws $addr $val
$addr = "G20"; $val = "0.0794"
Set-TextCell $ws $addr $val
$addr = "H20"; $val = 10
Set-NumCell $ws $addr $val
# row 21 / index 19
$addr = "A21"; $val = 19
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B21"; $val = "690008"
Set-TextCell $ws $addr $val
$addr = "C21"; $val = "民生加银中证内地资源主题指数"
Set-TextCell $ws $addr $val
$addr = "D21"; $val = "1.87"
Set-TextCell $ws $addr $val
$addr = "E21"; $val = "94.56"
Set-TextCell $ws $addr $val
$addr = "F21"; $val = "2.65"
Set-TextCell $ws $addr $val
$addr = "G21"; $val = "0.0496"
Set-TextCell $ws $addr $val
$addr = "H21"; $val = 10
Set-NumCell $ws $addr $val
# row 22 / index 20
$addr = "A22"; $val = 20
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B22"; $val = "159881"
Set-TextCell $ws $addr $val
$addr = "C22"; $val = "国泰中证有色金属交易型开放式指数证券投资基金"
Set-TextCell $ws $addr $val
$addr = "D22"; $val = "1.29"
Set-TextCell $ws $addr $val
$addr = "E22"; $val = "98.16"
Set-TextCell $ws $addr $val
$addr = "F22"; $val = "3.23"
Set-TextCell $ws $addr $val
$addr = "G22"; $val = "0.0417"
Set-TextCell $ws $addr $val
$addr = "H22"; $val = 6
Set-NumCell $ws $addr $val
# row 23 / index 21
$addr = "A23"; $val = 21
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B23"; $val = "007423"
Set-TextCell $ws $addr $val
$addr = "C23"; $val = "西部利得聚禾灵活配置混合A"
Set-TextCell $ws $addr $val
$addr = "D23"; $val = "0.60"
Set-TextCell $ws $addr $val
$addr = "E23"; $val = "69.21"
Set-TextCell $ws $addr $val
$addr = "F23"; $val = "3.78"
Set-TextCell $ws $addr $val
$addr = "G23"; $val = "0.0227"
Set-TextCell $ws $addr $val
$addr = "H23"; $val = 10
Set-NumCell $ws $addr $val
# row 24 / index 22
$addr = "A24"; $val = 22
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B24"; $val = "011027"
Set-TextCell $ws $addr $val
$addr = "C24"; $val = "国寿安保稳弘混合A"
Set-TextCell $ws $addr $val
$addr = "D24"; $val = "1.00"
Set-TextCell $ws $addr $val
$addr = "E24"; $val = "36.25"
Set-TextCell $ws $addr $val
$addr = "F24"; $val = "2.03"
Set-TextCell $ws $addr $val
$addr = "G24"; $val = "0.0203"
Set-TextCell $ws $addr $val
$addr = "H24"; $val = 9
Set-NumCell $ws $addr $val
# row 25 / index 23
$addr = "A25"; $val = 23
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B25"; $val = "000894"
Set-TextCell $ws $addr $val
$addr = "C25"; $val = "中欧睿达定期开放混合A"
Set-TextCell $ws $addr $val
$addr = "D25"; $val = "0.83"
Set-TextCell $ws $addr $val
$addr = "E25"; $val = "31.76"
Set-TextCell $ws $addr $val
$addr = "F25"; $val = "2.33"
Set-TextCell $ws $addr $val
$addr = "G25"; $val = "0.0193"
Set-TextCell $ws $addr $val
$addr = "H25"; $val = 4
Set-NumCell $ws $addr $val
# row 26 / index 24
$addr = "A26"; $val = 24
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B26"; $val = "014331"
Set-TextCell $ws $addr $val
$addr = "C26"; $val = "华泰柏瑞中证稀土产业ETF联接A"
Set-TextCell $ws $addr $val
$addr = "D26"; $val = "0.86"
Set-TextCell $ws $addr $val
$addr = "E26"; $val = "24.22"
Set-TextCell $ws $addr $val
$addr = "F26"; $val = "1.82"
Set-TextCell $ws $addr $val
$addr = "G26"; $val = "0.0157"
Set-TextCell $ws $addr $val
$addr = "H26"; $val = 2
Set-NumCell $ws $addr $val
# row 27 / index 25
$addr = "A27"; $val = 25
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B27"; $val = "007424"
Set-TextCell $ws $addr $val
$addr = "C27"; $val = "西部利得聚禾灵活配置混合C"
Set-TextCell $ws $addr $val
$addr = "D27"; $val = "0.41"
Set-TextCell $ws $addr $val
$addr = "E27"; $val = "69.21"
Set-TextCell $ws $addr $val
$addr = "F27"; $val = "3.78"
Set-TextCell $ws $addr $val
$addr = "G27"; $val = "0.0155"
Set-TextCell $ws $addr $val
$addr = "H27"; $val = 10
Set-NumCell $ws $addr $val
# row 28 / index 26
$addr = "A28"; $val = 26
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B28"; $val = "159876"
Set-TextCell $ws $addr $val
$addr = "C28"; $val = "华宝中证有色金属ETF"
Set-TextCell $ws $addr $val
$addr = "D28"; $val = "0.45"
Set-TextCell $ws $addr $val
$addr = "E28"; $val = "98.90"
Set-TextCell $ws $addr $val
$addr = "F28"; $val = "3.28"
Set-TextCell $ws $addr $val
$addr = "G28"; $val = "0.0148"
Set-TextCell $ws $addr $val
$addr = "H28"; $val = 6
Set-NumCell $ws $addr $val
# row 29 / index 27
$addr = "A29"; $val = 27
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B29"; $val = "159871"
Set-TextCell $ws $addr $val
$addr = "C29"; $val = "银华中证有色金属ETF"
Set-TextCell $ws $addr $val
$addr = "D29"; $val = "0.43"
Set-TextCell $ws $addr $val
$addr = "E29"; $val = "97.56"
Set-TextCell $ws $addr $val
$addr = "F29"; $val = "3.24"
Set-TextCell $ws $addr $val
$addr = "G29"; $val = "0.0139"
Set-TextCell $ws $addr $val
$addr = "H29"; $val = 6
Set-NumCell $ws $addr $val
# row 30 / index 28
$addr = "A30"; $val = 28
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B30"; $val = "011028"
Set-TextCell $ws $addr $val
$addr = "C30"; $val = "国寿安保稳弘混合C"
Set-TextCell $ws $addr $val
$addr = "D30"; $val = "0.68"
Set-TextCell $ws $addr $val
$addr = "E30"; $val = "36.25"
Set-TextCell $ws $addr $val
$addr = "F30"; $val = "2.03"
Set-TextCell $ws $addr $val
$addr = "G30"; $val = "0.0138"
Set-TextCell $ws $addr $val
$addr = "H30"; $val = 9
Set-NumCell $ws $addr $val
# row 31 / index 29
$addr = "A31"; $val = 29
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B31"; $val = "014332"
Set-TextCell $ws $addr $val
$addr = "C31"; $val = "华泰柏瑞中证稀土产业ETF联接C"
Set-TextCell $ws $addr $val
$addr = "D31"; $val = "0.70"
Set-TextCell $ws $addr $val
$addr = "E31"; $val = "24.22"
Set-TextCell $ws $addr $val
$addr = "F31"; $val = "1.82"
Set-TextCell $ws $addr $val
$addr = "G31"; $val = "0.0127"
Set-TextCell $ws $addr $val
$addr = "H31"; $val = 2
Set-NumCell $ws $addr $val
# row 32 / index 30
$addr = "A32"; $val = 30
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B32"; $val = "159880"
Set-TextCell $ws $addr $val
$addr = "C32"; $val = "鹏华国证有色金属行业ETF"
Set-TextCell $ws $addr $val
$addr = "D32"; $val = "0.33"
Set-TextCell $ws $addr $val
$addr = "E32"; $val = "96.25"
Set-TextCell $ws $addr $val
$addr = "F32"; $val = "3.58"
Set-TextCell $ws $addr $val
$addr = "G32"; $val = "0.0118"
Set-TextCell $ws $addr $val
$addr = "H32"; $val = 5
Set-NumCell $ws $addr $val
# row 33 / index 31
$addr = "A33"; $val = 31
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B33"; $val = "516680"
Set-TextCell $ws $addr $val
$addr = "C33"; $val = "建信中证细分有色金属产业主题ETF"
Set-TextCell $ws $addr $val
$addr = "D33"; $val = "0.32"
Set-TextCell $ws $addr $val
$addr = "E33"; $val = "96.13"
Set-TextCell $ws $addr $val
$addr = "F33"; $val = "3.29"
Set-TextCell $ws $addr $val
$addr = "G33"; $val = "0.0105"
Set-TextCell $ws $addr $val
$addr = "H33"; $val = 6
Set-NumCell $ws $addr $val
# row 34 / index 32
$addr = "A34"; $val = 32
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B34"; $val = "516650"
Set-TextCell $ws $addr $val
$addr = "C34"; $val = "华夏中证细分有色金属产业主题交易型开放式指数证券投资基金"
Set-TextCell $ws $addr $val
$addr = "D34"; $val = "0.29"
Set-TextCell $ws $addr $val
$addr = "E34"; $val = "98.99"
Set-TextCell $ws $addr $val
$addr = "F34"; $val = "3.44"
Set-TextCell $ws $addr $val
$addr = "G34"; $val = "0.0100"
Set-TextCell $ws $addr $val
$addr = "H34"; $val = 6
Set-NumCell $ws $addr $val
# row 35 / index 33
$addr = "A35"; $val = 33
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B35"; $val = "013437"
Set-TextCell $ws $addr $val
$addr = "C35"; $val = "财通资管中证有色金属指数A"
Set-TextCell $ws $addr $val
$addr = "D35"; $val = "0.21"
Set-TextCell $ws $addr $val
$addr = "E35"; $val = "92.75"
Set-TextCell $ws $addr $val
$addr = "F35"; $val = "3.32"
Set-TextCell $ws $addr $val
$addr = "G35"; $val = "0.0070"
Set-TextCell $ws $addr $val
$addr = "H35"; $val = 6
Set-NumCell $ws $addr $val
# row 36 / index 34
$addr = "A36"; $val = 34
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B36"; $val = "013438"
Set-TextCell $ws $addr $val
$addr = "C36"; $val = "财通资管中证有色金属指数C"
Set-TextCell $ws $addr $val
$addr = "D36"; $val = "0.10"
Set-TextCell $ws $addr $val
$addr = "E36"; $val = "92.75"
Set-TextCell $ws $addr $val
$addr = "F36"; $val = "3.32"
Set-TextCell $ws $addr $val
$addr = "G36"; $val = "0.0033"
Set-TextCell $ws $addr $val
$addr = "H36"; $val = 6
Set-NumCell $ws $addr $val
# row 37 / index 35
$addr = "A37"; $val = 35
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B37"; $val = "009648"
Set-TextCell $ws $addr $val
$addr = "C37"; $val = "中欧睿达定期开放混合C"
Set-TextCell $ws $addr $val
$addr = "D37"; $val = "0.00"
Set-TextCell $ws $addr $val
$addr = "E37"; $val = "31.76"
Set-TextCell $ws $addr $val
$addr = "F37"; $val = "2.33"
Set-TextCell $ws $addr $val
$addr = "G37"; $val = 0
Set-NumCell $ws $addr $val
$addr = "H37"; $val = 4
Set-NumCell $ws $addr $val
# row 38 / index 36
$addr = "A38"; $val = 36
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B38"; $val = "007234"
Set-TextCell $ws $addr $val
$addr = "C38"; $val = "博时优势企业3年封闭运作灵活配置混合（LOF）C"
Set-TextCell $ws $addr $val
$addr = "E38"; $val = "82.86"
Set-TextCell $ws $addr $val
$addr = "F38"; $val = "3.96"
Set-TextCell $ws $addr $val
$addr = "G38"; $val = 0
Set-NumCell $ws $addr $val
$addr = "H38"; $val = 9
Set-NumCell $ws $addr $val
# row 39 / index 37
$addr = "A39"; $val = 37
Set-StyledCell $ws $addr $val $styleSrc
$addr = "B39"; $val = "015407"
Set-TextCell $ws $addr $val
$addr = "C39"; $val = "国寿安保稳弘混合E"
Set-TextCell $ws $addr $val
$addr = "D39"; $val = "0.00"
Set-TextCell $ws $addr $val
$addr = "E39"; $val = "36.25"
Set-TextCell $ws $addr $val
$addr = "F39"; $val = "2.03"
Set-TextCell $ws $addr $val
$addr = "G39"; $val = 0
Set-NumCell $ws $addr $val
$addr = "H39"; $val = 9
Set-NumCell $ws $addr $val

# ----------------------------------------------------------------------
# Step 2 - insert a brand-new sheet named "总计" right after "2022-Q1" and
# rebuild the totals table: same 3 columns (日期 / 持有数量(只) /
# 持有市值(亿元)) as before, with one new leading row for 2022-Q1 and every
# other row shifted down by one index/row.
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "总计"

# -- header row --
$addr = "B1"; $val = "日期"
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "C1"; $val = "持有数量(只)"
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "D1"; $val = "持有市值(亿元)"
Set-StyledCell $ws2 $addr $val $styleSrc

# -- data rows (quarterly totals, newest first) --
# row 2 / index 0
$addr = "A2"; $val = 0
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "B2"; $val = "2022-Q1"
Set-TextCell $ws2 $addr $val
$addr = "C2"; $val = 38
Set-NumCell $ws2 $addr $val
$addr = "D2"; $val = 9.75
Set-NumCell $ws2 $addr $val
# row 3 / index 1
$addr = "A3"; $val = 1
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "B3"; $val = "2021-Q4"
Set-TextCell $ws2 $addr $val
$addr = "C3"; $val = 27
Set-NumCell $ws2 $addr $val
$addr = "D3"; $val = 6.6
Set-NumCell $ws2 $addr $val
# row 4 / index 2
$addr = "A4"; $val = 2
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "B4"; $val = "2021-Q3"
Set-TextCell $ws2 $addr $val
$addr = "C4"; $val = 64
Set-NumCell $ws2 $addr $val
$addr = "D4"; $val = 16.88
Set-NumCell $ws2 $addr $val
# row 5 / index 3
$addr = "A5"; $val = 3
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "B5"; $val = "2021-Q2"
Set-TextCell $ws2 $addr $val
$addr = "C5"; $val = 24
Set-NumCell $ws2 $addr $val
$addr = "D5"; $val = 5.49
Set-NumCell $ws2 $addr $val
# row 6 / index 4
$addr = "A6"; $val = 4
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "B6"; $val = "2021-Q1"
Set-TextCell $ws2 $addr $val
$addr = "C6"; $val = 6
Set-NumCell $ws2 $addr $val
$addr = "D6"; $val = 3.4
Set-NumCell $ws2 $addr $val
# row 7 / index 5
$addr = "A7"; $val = 5
Set-StyledCell $ws2 $addr $val $styleSrc
$addr = "B7"; $val = "2020-Q4"
Set-TextCell $ws2 $addr $val
$addr = "C7"; $val = 5
Set-NumCell $ws2 $addr $val
$addr = "D7"; $val = 2.48
Set-NumCell $ws2 $addr $val

Write-Output "edit complete"
